$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the example/demo data that had been filled into the risk table
# (the "sensor reading failure" sample row), turning those cells back
# into blank template cells while keeping their formatting intact.
$ws.Range("B6").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""

# Row 8's height had auto-grown to fit the wrapped example text; now that
# the text is gone, let it shrink back down to the default row height.
$ws.Rows(8).AutoFit()

# Leave the selection where the user ended up after clearing the data.
$ws.Range("J4:L4").Select()
